# Update the Shopkeepers sheet:
#  - Rows 2-5 get new ID/Name/Contact_Info/IsDeleted/Brand values
#  - Rows 6-9 are removed entirely (workbook shrinks from A1:E9 to A1:E5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force the cell to be stored as text (preserves leading zeros in
    # phone numbers) without leaving a permanent Text number format
    # applied to the cell once we are done.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2: Ahmed
$ws.Range("A2").Value = 12
Set-TextValue $ws.Range("B2") "Ahmed"
Set-TextValue $ws.Range("C2") "03448657309"
$ws.Range("D2").Value = 0
Set-TextValue $ws.Range("E2") "Bonapapa"

# Row 3: Gul Khan
$ws.Range("A3").Value = 13
Set-TextValue $ws.Range("B3") "Gul Khan"
Set-TextValue $ws.Range("C3") "23455534566"
$ws.Range("D3").Value = 0
Set-TextValue $ws.Range("E3") "Candyland"

# Row 4: Karim
$ws.Range("A4").Value = 14
Set-TextValue $ws.Range("B4") "Karim"
Set-TextValue $ws.Range("C4") "03465537715"
$ws.Range("D4").Value = 0
Set-TextValue $ws.Range("E4") "Candyland"

# Row 5: Nadim
$ws.Range("A5").Value = 15
Set-TextValue $ws.Range("B5") "Nadim"
Set-TextValue $ws.Range("C5") "03475858399"
$ws.Range("D5").Value = 0
Set-TextValue $ws.Range("E5") "Candyland"

# Remove the now-obsolete rows 6-9 entirely
$ws.Range("A6:E9").EntireRow.Delete()
